$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary pairs (English term / German translation), each inserted as a
# new row just above its alphabetically-following neighbour so that the
# existing (near-alphabetical) ordering of the list is preserved.
#
# New-sheet row numbers for the five additions: 1, 4, 17, 21, 30.

# 1) "ameliorate" / "ausbessern" -> becomes the new row 1 (goes before "asset")
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "ameliorate"
$ws.Range("B1").Value = "ausbessern"

# 2) "bailout" / "Notverkauf" -> becomes the new row 4 (goes before "bond")
$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "bailout"
$ws.Range("B4").Value = "Notverkauf"

# 3) "delinquencies" / "Nichtzahlung bei Fälligkeit" -> new row 17 (before "detrimental")
$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = "delinquencies"
$ws.Range("B17").Value = "Nichtzahlung bei Fälligkeit"

# 4) "equity" / "Eigenanteil" -> new row 21 (before "exposure")
$ws.Rows("21:21").Insert()
$ws.Range("A21").Value = "equity"
$ws.Range("B21").Value = "Eigenanteil"

# 5) "interest rate" / "Zinssatz" -> new row 30 (before "lend")
$ws.Rows("30:30").Insert()
$ws.Range("A30").Value = "interest rate"
$ws.Range("B30").Value = "Zinssatz"

# Hidden defined name added by the MySQL-for-Excel add-in metadata that was
# present in the saved workbook.
$formula = '=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&" "&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)'
$name = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $formula)
$name.Visible = $false

# Restore the selection to the cell that ends up holding the new "equity" row
# translation, matching the author's last-edited cell.
$ws.Range("B21").Select()
